$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 153 - this pushes the existing data rows
# 153..168 down to 154..169 (matching the rest of the diff, which is just
# every subsequent row's values shifting down by one position).
$ws.Rows.Item(153).Insert()

# Populate the newly inserted row 153 with a new data record. Columns
# A, B, C, E, F, G, H, N, O, Q, R repeat the same constant values used by
# every other row in this block (market/region/product metadata); only
# D (date), I (category), J, K, L, M, P (prices) differ for this record.
$ws.Cells.Item(153, 1).Value = 4
$ws.Cells.Item(153, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(153, 3).Value = "Los Lagos"
$ws.Cells.Item(153, 4).Value = 44449
$ws.Cells.Item(153, 5).Value = 10
$ws.Cells.Item(153, 6).Value = 100112045
$ws.Cells.Item(153, 7).Value = "Zapallo"
$ws.Cells.Item(153, 8).Value = "Paine"
$ws.Cells.Item(153, 9).Value = "1a (guarda)"
$ws.Cells.Item(153, 10).Value = 900
$ws.Cells.Item(153, 11).Value = 600
$ws.Cells.Item(153, 12).Value = 600
$ws.Cells.Item(153, 13).Value = 600
$ws.Cells.Item(153, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(153, 15).Value = "Región Metropolitana"
$ws.Cells.Item(153, 16).Value = 600
$ws.Cells.Item(153, 17).Value = 1
$ws.Cells.Item(153, 18).Value = "Hortaliza"
